$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (style index 0 / no special number format) used to restore
# cell formatting after writing date-shaped text, so Excel's automatic date
# recognition does not turn "YYYY-MM-DD" strings into real date serials.
$plainStyle = $ws.Range("A1").Style

# Row 2
$ws.Range("A2").Value = 80491706
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = '2019-09-16'
$ws.Range("AA2").Style = $plainStyle
$ws.Range("B2").Value = 89356
$ws.Range("D2").Value = 'LC'
$ws.Range("E2").Value = 5447
$ws.Range("F2").Value = 'Vedticka'
$ws.Range("G2").Value = 'Fuscoporia viticola'
$ws.Range("H2").Value = '(Schwein.) Murrill'
$ws.Range("P2").Value = 'Storåsen, Dlr'
$ws.Range("Q2").Value = 414384.947644942
$ws.Range("R2").Value = 6715092.921643512
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = '2019-09-16'
$ws.Range("Y2").Style = $plainStyle
# Row 3
$ws.Range("A3").Value = 80491716
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = '2019-09-16'
$ws.Range("AA3").Style = $plainStyle
$ws.Range("AC3").ClearContents()
$ws.Range("B3").Value = 77506
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = 'Garnlav'
$ws.Range("G3").Value = 'Alectoria sarmentosa'
$ws.Range("H3").Value = '(Ach.) Ach.'
$ws.Range("P3").Value = 'Storåsen, Dlr'
$ws.Range("Q3").Value = 414430.7598552333
$ws.Range("R3").Value = 6715154.798017079
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = '2019-09-16'
$ws.Range("Y3").Style = $plainStyle
# Row 4
$ws.Range("A4").Value = 80491714
$ws.Range("Q4").Value = 414472.0158322724
$ws.Range("R4").Value = 6715170.051899989
# Row 5
$ws.Range("A5").Value = 80491721
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = '2019-09-17'
$ws.Range("AA5").Style = $plainStyle
$ws.Range("P5").Value = 'Gårdtjärnsmyrorna, Dlr'
$ws.Range("Q5").Value = 414915.9210306811
$ws.Range("R5").Value = 6715399.095652443
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = '2019-09-17'
$ws.Range("Y5").Style = $plainStyle
# Row 6
$ws.Range("A6").Value = 80491713
$ws.Range("B6").Value = 89392
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 1202
$ws.Range("F6").Value = 'Ullticka'
$ws.Range("G6").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H6").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q6").Value = 414475.0636448057
$ws.Range("R6").Value = 6715173.91479082
# Row 7
$ws.Range("A7").Value = 80491722
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = '2019-09-17'
$ws.Range("AA7").Style = $plainStyle
$ws.Range("AC7").Value = 'På död gran'
$ws.Range("B7").Value = 85703
$ws.Range("E7").Value = 510
$ws.Range("F7").Value = 'Doftskinn'
$ws.Range("G7").Value = 'Cystostereum murrayi'
$ws.Range("H7").Value = '(Berk. & M.A. Curtis.) Pouzar'
$ws.Range("P7").Value = 'Gårdtjärnsmyrorna, Dlr'
$ws.Range("Q7").Value = 414916.8826095874
$ws.Range("R7").Value = 6715398.089234225
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = '2019-09-17'
$ws.Range("Y7").Style = $plainStyle
